$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Excel alignment constants (xlHAlignLeft / xlHAlignRight / xlHAlignCenter)
$xlLeft = -4131

# --- Update the parameter-statistics labels (column A, rows 8-11) ---
$ws.Range("A8").Value = "Range"
$ws.Range("A9").Value = "Mean"
$ws.Range("A10").Value = "Standard Deviation"
$ws.Range("A11").Value = "Mean rate of change"

# --- Update the matching placeholder labels (column B, rows 8-11) ---
$ws.Range("B8").Value = "B8 range"
$ws.Range("B9").Value = "B9 mean"
$ws.Range("B10").Value = " B10 std derv"
$ws.Range("B11").Value = "B11 mean rate of change"

# B8:B11 keep their box border but now also left-align their text
$ws.Range("B8:B11").HorizontalAlignment = $xlLeft

# Header row (A13:C13) and the templated data row (A14:C14) switch from
# center-aligned to left-aligned text
$ws.Range("A13:C13").HorizontalAlignment = $xlLeft
$ws.Range("A14:C14").HorizontalAlignment = $xlLeft

# --- Adjust the window view so it mirrors where the user scrolled/selected ---
$win = $excel.ActiveWindow
$win.ScrollRow = 7
$win.ScrollColumn = 1
$ws.Range("C17").Select()
